$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.883.64'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '3.153.97'
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = "'216.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").Value = "'626.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.30%  '

$ws.Range("D7").Value = "'1.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +32.60%  '

$ws.Range("D8").Value = "'0.368"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.72%  '

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").Value = '3.152.54'
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("D11").Value = "'0.765"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +15.14%  '

$ws.Range("E12").Value = '  +6.89%  '

$ws.Range("D13").Value = "'5.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.44%  '

$ws.Range("E14").Value = '  -1.22%  '

$ws.Range("D15").Value = "'35.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.40%  '

$ws.Range("D16").Value = '90.690.49'
$ws.Range("E16").Value = '  -0.21%  '

$ws.Range("D17").Value = '3.740.03'
$ws.Range("E17").Value = '  +2.42%  '

$ws.Range("D18").Value = '3.203.05'
$ws.Range("E18").Value = '  +4.05%  '

$ws.Range("D19").Value = "'3.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.53%  '

$ws.Range("D20").Value = "'14.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.93%  '

$ws.Range("D21").Value = "'475.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.45%  '

$ws.Range("E22").Value = '  -4.31%  '

$ws.Range("D23").Value = "'9.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.22%  '

$ws.Range("D24").Value = "'5.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.28%  '

$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").Value = "'5.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.89%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = "'94.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.83%  '

$ws.Range("E27").Value = '  +4.62%  '

$ws.Range("D28").Value = '3.326.67'
$ws.Range("E28").Value = '  +2.87%  '

$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("D30").Value = "'9.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.47%  '

$ws.Range("D31").Value = "'0.163"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.10%  '

$ws.Range("E32").Value = '  +53.56%  '

$ws.Range("E33").Value = '  -7.41%  '

$ws.Range("D34").Value = "'27.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +19.53%  '

$ws.Range("D35").Value = "'519.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("E36").Value = '  +6.22%  '

$ws.Range("E37").Value = '  +5.86%  '

$ws.Range("D38").Value = "'3.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.53%  '

$ws.Range("D39").Value = "'6.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.91%  '

$ws.Range("D40").Value = "'1.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.39%  '

$ws.Range("D41").Value = "'0.0912"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +27.35%  '

$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").Value = "'0.427"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.09%  '

$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").Value = "'22.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("E44").Value = '  -0.15%  '

$ws.Range("D45").Value = "'1.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.28%  '

$ws.Range("D46").Value = "'0.738"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +21.62%  '

$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").Value = "'4.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.84%  '

$ws.Range("D49").Value = "'150.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.87%  '

$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").Value = "'1.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.67%  '

$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = "'45.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.18%  '
